$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Updated query text for the "CasesTab" row (row 2, column B) - adds Age /
# Weight numeric handling and a Cohort column.
# ---------------------------------------------------------------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
   WHERE f.file_type IN ["Pathology Report"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

# ---------------------------------------------------------------------------
# Shared "StatQuery" text (column C) - identical for every tab row.
# ---------------------------------------------------------------------------
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(f)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE  f.file_type IN ["Pathology Report"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# ---------------------------------------------------------------------------
# Updated query text for the "SamplesTab" row (row 3, column B).
# ---------------------------------------------------------------------------
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
MATCH (f:file)-[*]->(c)
WHERE f.file_type IN ["Pathology Report"]   
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

# ---------------------------------------------------------------------------
# Updated query text for the "FilesTab" row (row 4, column B) - now computes
# a human readable file size and adds a Sample ID column.
# ---------------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-[*]->(c:case)
MATCH (f)-->(parent)
OPTIONAL MATCH (f)-->(samp:sample)
OPTIONAL MATCH (p:program)<--(s:study)<--(c)
OPTIONAL MATCH (c)<--(demo:demographic)
OPTIONAL MATCH (c)<--(diag:diagnosis)
WITH
    f, c, parent, samp, p, s, demo, diag
WHERE f.file_type IN ["Pathology Report"] 
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    2 as precision
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag, unit,
    round(factor * value)/factor AS size
RETURN
    coalesce(f.file_name, '') AS `File Name`,
    coalesce(f.file_format, '') AS `Format`,
    coalesce(f.file_type, '') AS `File Type`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    coalesce(labels(parent)[0], '') AS `Association`,
    coalesce(f.file_description, '') AS `Description`,
    coalesce(samp.sample_id, '') AS `Sample ID`,
    coalesce(c.case_id, '') AS `Case ID`,
    coalesce(demo.breed,'') AS Breed ,
    coalesce(diag.disease_term,'') AS Diagnosis
    ORDER BY f.file_name asc
     limit 100
'@

# ---------------------------------------------------------------------------
# Brand new query text for the new "StudyFilesTab" row (row 5, column B).
# ---------------------------------------------------------------------------
$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (samp:sample)-->(c)
MATCH (c)<--(demo:demographic)
WHERE f.file_type IN ["Pathology Report"]
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

# ---------------------------------------------------------------------------
# Write the query text into the sheet (column B = per-tab query,
# column C = shared StatQuery text for every tab).
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# New row 5 - "StudyFilesTab"
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = "TC03_Canine_Filter_FileType-PathologyRep_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC03_Canine_Filter_FileType-PathologyRep_WebData.xlsx"

# ---------------------------------------------------------------------------
# Formatting: bump the base font size for the whole used range to 15pt, and
# make sure the long query cells (column B and C) wrap their text.
# ---------------------------------------------------------------------------
$ws.Range("A1:E5").Font.Size = 15
$ws.Range("B2:C5").WrapText = $true

$ws.Rows(2).RowHeight = 409.5
$ws.Rows(3).RowHeight = 409.5
$ws.Rows(4).RowHeight = 409.5
$ws.Rows(5).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Sheet-level view tweaks to match the authored edit.
# ---------------------------------------------------------------------------
$ws.StandardHeight = 19.5
$ws.Range("D12").Select()
